# Add new columns I ("I0") and J ("IF") to the active worksheet,
# mirroring the style of the existing header row and filling in the
# per-row data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, centered, bordered) from the existing
# header cell H1 so the new headers match the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows (rows 2-38) ---------------------------------------------
$data = @{
    2  = @(1, 4)
    3  = @(11, 11)
    4  = @(9, 9)
    5  = @(8, 8)
    6  = @(7, 8)
    7  = @(7, 8)
    8  = @(9, 9)
    9  = @(9, 9)
    10 = @(11, 11)
    11 = @(6, 9)
    12 = @(9, 9)
    13 = @(7, 8)
    14 = @(1, 6)
    15 = @(1, 6)
    16 = @(1, 6)
    17 = @(1, 4)
    18 = @(1, 6)
    19 = @(1, 7)
    20 = @(1, 6)
    21 = @(1, 6)
    22 = @(1, 6)
    23 = @(1, 5)
    24 = @(1, 7)
    25 = @(1, 5)
    26 = @(1, 7)
    27 = @(1, 5)
    28 = @(1, 6)
    29 = @(1, 5)
    30 = @(1, 5)
    31 = @(1, 7)
    32 = @(1, 7)
    33 = @(1, 6)
    34 = @(1, 5)
    35 = @(1, 4)
    36 = @(1, 5)
    37 = @(1, 4)
    38 = @(1, 3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]   # column I
    $ws.Cells.Item($row, 10).Value = $vals[1]  # column J
}

$wb.Save()
